$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new paragraph "-Mode sombre/blanc" immediately before the
#    existing "-Convertisseur" paragraph (so "-Convertisseur" becomes the
#    second bullet, preceded by the new one).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("-Convertisseur", $false, $false, $false, $false, $false,
                         $true, 1, $false, "-Mode sombre/blanc`r-Convertisseur", 2)

# ---------------------------------------------------------------------------
# 2) Collapse the run-split / proofErr-wrapped currency list
#    (" (Bitcoin (BTC), Euros, Dollars, " + spell-checked "Tether" +
#    " (USDT), " + spell-checked "Ethereum" + " (ETH), USDC, BUSD)") into a
#    single plain run, while leaving the preceding
#    "-Cours de plusieurs monnaies en direct" run untouched/separate.
#
#    Find/Replace across several runs in this host tends to also gobble up
#    the unrelated, immediately-preceding run when it shares the same
#    (empty) formatting, so instead of replacing in place we build the
#    clean text in a fresh paragraph and then splice it back in:
# ---------------------------------------------------------------------------

# 2a) Append a new paragraph right after the "-Cours..." paragraph holding
#     the fully merged, proofErr-free replacement text as a single run.
$r = $d.Content
$r.Find.Execute(" (ETH), USDC, BUSD)", $false, $false, $false, $false, $false,
                $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.MoveStart(1, 1)
$r.InsertAfter(" (Bitcoin (BTC), Euros, Dollars, Tether (USDT), Ethereum (ETH), USDC, BUSD)")

# 2b) Delete the old proofErr-laden run content from the original paragraph,
#     leaving just the "-Cours de plusieurs monnaies en direct" run behind.
$r2 = $d.Content
$r2.Find.Execute(" (Bitcoin (BTC), Euros, Dollars, ", $false, $false, $false, $false, $false,
                 $true, 1, $false, "", 0)
$start = $r2.Start
$r3 = $d.Content
$r3.Find.Execute(" (ETH), USDC, BUSD)", $false, $false, $false, $false, $false,
                 $true, 1, $false, "", 0)
$end = $r3.End
$d.Range($start, $end).Delete()

# 2c) Re-join the (now trimmed) original paragraph with the new clean
#     paragraph by deleting the paragraph mark between them, yielding one
#     paragraph with exactly two runs.
$r4 = $d.Content
$r4.Find.Execute("-Cours de plusieurs monnaies en direct", $false, $false, $false, $false, $false,
                 $true, 1, $false, "", 0)
$markStart = $r4.End
$d.Range($markStart, $markStart + 1).Delete()
